# Aula 34 - Fragmentando paginas de cadastros e listagens
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D23 gains wrap-text formatting (style index 1 in the original workbook)
$ws.Range("D23").WrapText = $true

# New row 24: Aula 34
$ws.Range("B24").Value = 34
$ws.Range("C24").Value = "7. Thymeleaf para as Views"
$ws.Range("D24").Value = "34. Fragmentando páginas de cadastros e listagens"
$ws.Range("E24").Value = "demonstra na pratica como fragmentar as demais páginas usando o layout:fragment e layout:decorate"
$ws.Range("D24").WrapText = $true
$ws.Range("E24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 30

# New row 25: RESUMO
$ws.Range("B25").Value = 35
$ws.Range("C25").Value = "7. Thymeleaf para as Views"
$ws.Range("D25").Value = "RESUMO"
$ws.Range("E25").Value = "Um resumão sobre tudo que foi aprendido na sessão/módulo/capitulo`na leitura é valida caso surja dúvidas"
$ws.Range("D25").WrapText = $true
$ws.Range("E25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 45

$ws.Range("E25").Select()
